# Fruta / hortaliza, semanal
# Insert a new weekly record at row 41 ("Vega Monumental Concepción" / "Poroto
# granado" data block), pushing the existing rows 41-61 down to 42-62.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 41; this shifts rows 41:61 to 42:62
# and carries the existing formatting (incl. the date style on column D) down
# with them, leaving a blank row 41 with the same row formatting.
$ws.Rows.Item(41).Insert()

# Columns A, B, C, E, F, G, H, I, N, Q, R are identical for every record in
# this block, so copy them from the row that was just pushed down (row 42,
# the former row 41) into the new row 41.
$ws.Range("A41").Value = $ws.Range("A42").Value()
$ws.Range("B41").Value = $ws.Range("B42").Value()
$ws.Range("C41").Value = $ws.Range("C42").Value()
$ws.Range("E41").Value = $ws.Range("E42").Value()
$ws.Range("F41").Value = $ws.Range("F42").Value()
$ws.Range("G41").Value = $ws.Range("G42").Value()
$ws.Range("H41").Value = $ws.Range("H42").Value()
$ws.Range("I41").Value = $ws.Range("I42").Value()
$ws.Range("N41").Value = $ws.Range("N42").Value()
$ws.Range("Q41").Value = $ws.Range("Q42").Value()
$ws.Range("R41").Value = $ws.Range("R42").Value()

# New record values for the inserted row.
$ws.Range("D41").Value = 45001
$ws.Range("J41").Value = 140
$ws.Range("K41").Value = 30000
$ws.Range("L41").Value = 32000
$ws.Range("M41").Value = 30857
$ws.Range("O41").Value = "Región Metropolitana"
$ws.Range("P41").Value = 1234
